$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "Derek Bays" -> "Derek J Bays" (middle initial added to author list)
$ws.Range("A2").Value = "Amy P. Hsu, Agnieszka Korzeniowska, Cynthia C Aguilar, Jiande Gu, Eric Karlins, Andrew J. Oler, Gang Chen, Glennys V. Reynoso, Joie Davis, Alexandria Laurel Chaput, Tao Peng, Ling Sun, Justin Lack, Derek J Bays, Ethan R. Stewart, Sarah Waldman, Daniel A. Powell, Fariba M. Donovan, Jigar V. Desai, Nima Pouladi, Debra A. Long Priel, Daisuke Yamanaka, Sergio D. Rosenzweig, Julie E. Niemela, Jennifer Stoddard, Alexandra F. Freeman, Christa S. Zerbe, Douglas B. Kuhns, Yves A. Lussier, Kenneth N. Olivier, Richard C. Boucher, Heather D. Hickman, Jeffrey A. Frelinger, Joshua Fierer, Lisa F. Shubitz, Thomas L. Leto, George R. Thompson, John N. Galgiani, Michail S Lionakis, Steven M. Holland"

# Row 2: affiliation list was blank placeholders; now filled in with full text
$ws.Range("B2").Value = "Department of Cell Biology and Molecular Genetics, University of Maryland, College Park, Maryland, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; ; ; ; Marsico Lung Institute and Cystic Fibrosis Research Center, University of North Carolina at Chapel Hill, Chapel Hill, North Carolina, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Valley Fever Center for Excellence, University of Arizona College of Medicine–Tucson, Tucson, Arizona, USA.; Valley Fever Center for Excellence, University of Arizona College of Medicine–Tucson, Tucson, Arizona, USA.; Department of Respiratory and Critical Care Medicine, Laboratory of Pulmonary Immunology and Inflammation, West China Hospital, Sichuan University, Chengdu, Sichuan Province, China.; Marsico Lung Institute and Cystic Fibrosis Research Center, University of North Carolina at Chapel Hill, Chapel Hill, North Carolina, USA.; Advanced Biomedical Computational Science, Frederick National Laboratory for Cancer Research, Leidos Biomedical Research, Inc., Frederick, Maryland, USA.; NIAID Collaborative Bioinformatics Resource, NIAID, NIH, Bethesda, Maryland, USA.; Department of Internal Medicine, Division of Infectious Diseases, UC Davis Health, Sacramento, California, USA.; Department of Internal Medicine, Division of Infectious Diseases, UC Davis Health, Sacramento, California, USA.; Department of Internal Medicine, Division of Infectious Diseases, UC Davis Health, Sacramento, California, USA.; Department of Immunobiology, University of Arizona, Tucson, Arizona, USA.; Valley Fever Center for Excellence, University of Arizona College of Medicine–Tucson, Tucson, Arizona, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Center for Biomedical Informatics and Biostatistics and; The Center for Applied Genetics and Genomic Medicine, Department of Medicine, University of Arizona, Tucson, Arizona, USA.; Neutrophil Monitoring Laboratory, Applied/Developmental Research Directorate, Leidos Biomedical Research, Inc, Frederick National Laboratory for Cancer Research, Frederick, Maryland, USA.; Laboratory for Immunopharmacology of Microbial Products, School of Pharmacy, Tokyo University of Pharmacy and Life Sciences, Hachioji, Tokyo, Japan.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; ; ; ; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Neutrophil Monitoring Laboratory, Applied/Developmental Research Directorate, Leidos Biomedical Research, Inc, Frederick National Laboratory for Cancer Research, Frederick, Maryland, USA.; Center for Biomedical Informatics and Biostatistics and; The Center for Applied Genetics and Genomic Medicine, Department of Medicine, University of Arizona, Tucson, Arizona, USA.; Laboratory of Chronic Airway Infection, Pulmonary Branch, National Heart, Lung, and Blood Institute, NIH, Bethesda, Maryland, USA.; Marsico Lung Institute and Cystic Fibrosis Research Center, University of North Carolina at Chapel Hill, Chapel Hill, North Carolina, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Department of Immunobiology, University of Arizona, Tucson, Arizona, USA.; Division of Infectious Diseases, Departments of Pathology and Medicine, School of Medicine, University of California San Diego, La Jolla, California, USA.; VA HealthCare San Diego, San Diego, California, USA.; Valley Fever Center for Excellence, University of Arizona College of Medicine–Tucson, Tucson, Arizona, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Department of Medical Microbiology and Immunology, University of California Davis, Davis, California, USA.; Valley Fever Center for Excellence, University of Arizona College of Medicine–Tucson, Tucson, Arizona, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA.; Laboratory of Clinical Immunology and Microbiology, National Institute of Allergy and Infectious Diseases (NIAID), NIH, Bethesda, Maryland, USA."

# Row 2: cited_by_count 17 -> 19 (stored as text, matching the column's existing format)
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "19"

# Row 6: cited_by_count 1 -> 2 (stored as text, matching the column's existing format)
$ws.Range("M6").NumberFormat = "@"
$ws.Range("M6").Value = "2"
